{"js": "// \"show report s2 series correction\"\n// In the first table (Client / Enquiry / Project header block), the\n// placeholder \"1\" values for Client, Enquiry and Project become \"a\",\n// and the Enquiry \"Date\" value is corrected to a new date/time.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nasync function replaceCellText(rowIndex, colIndex, oldText, newText) {\n  const cell = table.getCell(rowIndex, colIndex);\n  const results = cell.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    // Found the old text as its own run \u2014 replace in place so the\n    // existing run formatting (bold, size, ...) is preserved.\n    results.items[0].insertText(newText, Word.InsertLocation.replace);\n  } else {\n    // Fallback: just overwrite the whole cell body.\n    cell.body.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// Row 0: Client | 1 -> a | Version | 5.1.2.0\nawait replaceCellText(0, 1, \"1\", \"a\");\n\n// Row 1: Enquiry | 1 -> a | Date | 12/20/2019, 04:33 PM -> 12/25/2019, 11:09 AM\nawait replaceCellText(1, 1, \"1\", \"a\");\nawait replaceCellText(1, 3, \"12/20/2019, 04:33 PM\", \"12/25/2019, 11:09 AM\");\n\n// Row 2: Project | 1 -> a | Model | TAC S2 C3\nawait replaceCellText(2, 1, \"1\", \"a\");\n", "ps1": "# \"show report s2 series correction\"\n# In the first table (Client / Enquiry / Project header block), the\n# placeholder \"1\" values for Client, Enquiry and Project become \"a\",\n# and the Enquiry \"Date\" value is corrected to a new date/time.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nfunction Set-CellText($table, $row, $col, $newText) {\n    $rng = $table.Cell($row, $col).Range\n    # Trim the trailing end-of-cell marker so assigning .Text replaces\n    # only the visible content and keeps the run's existing formatting\n    # (bold, size, ...) instead of resetting it.\n    $rng.MoveEnd(1, -1) | Out-Null\n    $rng.Text = $newText\n}\n\n# Row 1: Client | 1 -> a | Version | 5.1.2.0\nSet-CellText $t 1 2 \"a\"\n\n# Row 2: Enquiry | 1 -> a | Date | 12/20/2019, 04:33 PM -> 12/25/2019, 11:09 AM\nSet-CellText $t 2 2 \"a\"\nSet-CellText $t 2 4 \"12/25/2019, 11:09 AM\"\n\n# Row 3: Project | 1 -> a | Model | TAC S2 C3\nSet-CellText $t 3 2 \"a\"\n"}
